$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '68.169.43'
$ws.Range("E2").Value = '  +1.10%  '

# Row 3
$ws.Range("D3").Value = '3.344.09'
$ws.Range("E3").Value = '  +0.64%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.12%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.80%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.26'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.25%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.05%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.591'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.68%  '

# Row 9
$ws.Range("E9").Value = '  +4.00%  '

# Row 10
$ws.Range("E10").Value = '  +1.41%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '48.03'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.14%  '

# Row 12
$ws.Range("E12").Value = '  +1.76%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '693.58'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.96%  '

# Row 14
$ws.Range("D14").Value = '3.892.36'
$ws.Range("E14").Value = '  +0.72%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.42'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.65%  '

# Row 16
$ws.Range("D16").Value = '68.231.72'
$ws.Range("E16").Value = '  +0.96%  '

# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.445.59'
$ws.Range("E17").Value = '  +3.54%  '

# Row 18
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.119'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.35%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.43'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.31%  '

# Row 20
$ws.Range("E20").Value = '  +2.60%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.895'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.98%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.46'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.41%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.33%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '100.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.76%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.91'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.47%  '

# Row 26
$ws.Range("E26").Value = '  +1.41%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.51'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.34%  '

# Row 28
$ws.Range("E28").Value = '  -2.45%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.49'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.31%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.93'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.54%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '564.58'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.67%  '

# Row 32
$ws.Range("E32").Value = '  +1.43%  '

# Row 33
$ws.Range("E33").Value = '  +1.56%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '57.47'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.29%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.01%  '

# Row 36
$ws.Range("D36").Value = '3.687.81'
$ws.Range("E36").Value = '  +0.00%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.28'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.35%  '

# Row 38
$ws.Range("E38").Value = '  +4.32%  '

# Row 39
$ws.Range("E39").Value = '  +5.34%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.17'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.57%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.61'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.19%  '

# Row 42
$ws.Range("B42").Value = 'PEPE'
$ws.Range("C42").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D42").Value = '0.0₃0671'
$ws.Range("E42").Value = '  +1.89%  '

# Row 43
$ws.Range("B43").Value = 'TheGraph'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.335'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.97%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.26'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.21%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0414'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.40%  '

# Row 46
$ws.Range("E46").Value = '  +2.71%  '

# Row 47
$ws.Range("E47").Value = '  +1.02%  '

# Row 48
$ws.Range("E48").Value = '  -0.30%  '

# Row 49
$ws.Range("E49").Value = '  -0.09%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '130.98'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.34%  '

# Row 51
$ws.Range("E51").Value = '  +1.56%  '
